# Removed Test Case Inter-Dependency
# Update the product name / short name values on the input sheet so that
# the test data is unique (not shared with other test cases), and mirror
# the product name change onto the output sheet. Also switch the
# "active"/selected sheet from the input sheet to the output sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Order matters for how new shared strings get appended, so set the
# short name before the product name.
$ws1.Range("B2").Value = "245e"
$ws1.Range("B1").Value = "2450-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"
$ws1.Range("B13").Value = "Equal principal payments"

$ws2.Range("B1").Value = "2450-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"

# Move the selection/active cell and active sheet from the input sheet to
# the output sheet.
$null = $ws1.Range("B1").Select()
$null = $ws2.Activate()
